$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (autogluon) - previously empty inline strings, now populated
$ws.Range("B3").Value = "0.272 (0.244 ± 0.013)"
$ws.Range("C3").Value = "00:04:52 (00:05:01 ± 00:00:03)"
$ws.Range("D3").Value = "00:00:00 (00:00:01 ± 00:00:00)"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "31"

# Row 4 (autokeras) - fix mojibake "Â±" -> "±"
$ws.Range("B4").Value = "0.248 (0.234 ± 0.009)"
$ws.Range("C4").Value = "00:00:52 (00:01:07 ± 00:00:15)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"

# Row 6 (autosklearn) - fix mojibake "Â±" -> "±"
$ws.Range("B6").Value = "0.630 (0.597 ± 0.019)"
$ws.Range("C6").Value = "00:04:54 (00:05:01 ± 00:00:03)"
$ws.Range("D6").Value = "00:00:00 (00:00:05 ± 00:00:06)"

# Row 8 (fedot) - fix mojibake "Â±" -> "±"
$ws.Range("B8").Value = "0.576 (0.554 ± 0.014)"
$ws.Range("C8").Value = "00:05:08 (00:06:51 ± 00:01:15)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
